$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("GLOBAL RESULTS")
$ws.Range("C6").Value = 60132.19372407274
$ws.Range("C7").Value = 60925.74372407276
$ws.Range("C8").Value = 54833.16935166548
$ws.Range("C12").Value = 48966.6770447921
$ws.Range("C13").Value = 48966.6770447921
$ws.Range("C14").Value = 36096.6770447921
$ws.Range("C15").Value = 35367.5897537921
$ws.Range("C16").Value = 34517.379753792105
$ws.Range("C20").Value = 589695.3775841778
$ws.Range("C21").Value = 597477.4446916779
$ws.Range("C22").Value = 537729.7002225101
$ws.Range("C26").Value = 480199.06344131037
$ws.Range("C27").Value = 480199.06344131037
$ws.Range("C28").Value = 353987.4779413104
$ws.Range("C29").Value = 346837.57405902515
$ws.Range("C30").Value = 338499.8621625253

$ws = $wb.Worksheets.Item("FUSELAGE")
$ws.Range("C6").Value = 8026.0
$ws.Range("D6").Value = 60.3692529022719
$ws.Range("C7").Value = 5994.0
$ws.Range("D7").Value = 19.767418626491125
$ws.Range("C8").Value = 6125.0
$ws.Range("D8").Value = 22.38495813934904
$ws.Range("C9").Value = 6929.0
$ws.Range("D9").Value = 38.449857134293794
$ws.Range("C12").Value = 6682.666666666666
$ws.Range("D12").Value = 33.52781718517927

$ws = $wb.Worksheets.Item("WING")
$ws.Range("C7").Value = 5673.0
$ws.Range("D7").Value = 76.77853603814157
$ws.Range("C8").Value = 4398.0
$ws.Range("D8").Value = 37.04777040291673
$ws.Range("C9").Value = 5802.0
$ws.Range("D9").Value = 80.79835467888196
$ws.Range("C11").Value = 6078.0
$ws.Range("D11").Value = 89.3988968869777
$ws.Range("C12").Value = 5172.0
$ws.Range("D12").Value = 61.1666822473591
$ws.Range("C13").Value = 4618.714285714284
$ws.Range("D13").Value = 43.92553319355225

$ws = $wb.Worksheets.Item("HORIZONTAL TAIL")
$ws.Range("C8").Value = 237.0
$ws.Range("D8").Value = -56.77548787160313
$ws.Range("C9").Value = 484.0
$ws.Range("D9").Value = -11.72715666605872
$ws.Range("C10").Value = 483.33333333333326
$ws.Range("D10").Value = -11.848744604535218

$ws = $wb.Worksheets.Item("VERTICAL TAIL")
$ws.Range("C8").Value = 417.0
$ws.Range("D8").Value = -7.886017229953598
$ws.Range("C9").Value = 478.0
$ws.Range("D9").Value = 5.588690081731833

$ws = $wb.Worksheets.Item("NACELLES")
$ws.Range("C3").Value = 848.6666666666665
$ws.Range("D3").Value = -78.35254905961975
$ws.Range("C10").Value = 558.0
$ws.Range("D10").Value = -14.600550964187308
$ws.Range("C12").Value = 424.33333333333326
$ws.Range("C17").Value = 558.0
$ws.Range("D17").Value = -14.600550964187308
$ws.Range("C19").Value = 424.33333333333326

$ws = $wb.Worksheets.Item("LANDING GEARS")
$ws.Range("C5").Value = 1707.0
$ws.Range("D5").Value = -23.55917782454881
$ws.Range("C6").Value = 2397.0
$ws.Range("D6").Value = 7.339572791187165
$ws.Range("C7").Value = 2746.0
$ws.Range("D7").Value = 22.968071291030437
$ws.Range("C8").Value = 2375.0
$ws.Range("D8").Value = 6.3543952353231195
$ws.Range("C9").Value = 2306.25
$ws.Range("D9").Value = 3.2757153732479547
